$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript digits used in a couple of the tiny-price coins (e.g. 0.0<sub>3</sub>0904)
$sub3 = [char]0x2083
$sub6 = [char]0x2086

# Row -> [D value (or $null if unchanged), E value, D needs text-number-format?]
$updates = @(
    @{Row=2;  D='71.208.62';          E='  +2.50%  ';  Num=$false}
    @{Row=3;  D='4.007.90';           E='  +1.56%  ';  Num=$false}
    @{Row=4;  D=$null;                E='  -0.09%  ';  Num=$false}
    @{Row=5;  D='530.46';             E='  +4.39%  ';  Num=$true}
    @{Row=6;  D='148.98';             E='  +0.82%  ';  Num=$true}
    @{Row=7;  D='0.624';              E='  -0.40%  ';  Num=$true}
    @{Row=8;  D='1.00';               E='  +0.16%  ';  Num=$true}
    @{Row=9;  D=$null;                E='  +0.24%  ';  Num=$false}
    @{Row=10; D=$null;                E='  +0.02%  ';  Num=$false}
    @{Row=11; D='0.0000344';          E='  -1.74%  ';  Num=$true}
    @{Row=12; D='43.67';              E='  -0.04%  ';  Num=$true}
    @{Row=13; D='10.66';              E='  +1.06%  ';  Num=$true}
    @{Row=14; D='4.640.19';           E='  +1.48%  ';  Num=$false}
    @{Row=15; D='4.021.78';           E='  +1.61%  ';  Num=$false}
    @{Row=16; D='21.42';              E='  +7.01%  ';  Num=$true}
    @{Row=17; D='14.36';              E='  +0.63%  ';  Num=$true}
    @{Row=18; D='1.22';               E='  +0.71%  ';  Num=$true}
    @{Row=19; D='0.133';              E='  -1.95%  ';  Num=$true}
    @{Row=20; D='71.203.16';          E='  +2.43%  ';  Num=$false}
    @{Row=21; D='441.02';             E='  +1.11%  ';  Num=$true}
    @{Row=22; D=$null;                E='  +2.59%  ';  Num=$false}
    @{Row=23; D='93.02';              E='  +4.35%  ';  Num=$true}
    @{Row=24; D=$null;                E='  +4.39%  ';  Num=$false}
    @{Row=25; D='14.36';              E='  -2.44%  ';  Num=$true}
    @{Row=26; D=$null;                E='  +5.91%  ';  Num=$false}
    @{Row=27; D='10.90';              E='  -2.91%  ';  Num=$true}
    @{Row=28; D='37.01';              E='  -0.57%  ';  Num=$true}
    @{Row=29; D='13.65';              E='  +1.46%  ';  Num=$true}
    @{Row=30; D='694.49';             E='  -1.79%  ';  Num=$true}
    @{Row=31; D=$null;                E='  +0.16%  ';  Num=$false}
    @{Row=32; D=$null;                E='  -0.16%  ';  Num=$false}
    @{Row=33; D='6.86';               E='  +12.65%  '; Num=$true}
    @{Row=34; D='66.90';              E='  +1.94%  ';  Num=$true}
    @{Row=35; D="0.0${sub3}0904";     E='  +2.59%  ';  Num=$false}
    @{Row=36; D='0.443';              E='  -1.51%  ';  Num=$true}
    @{Row=37; D='41.04';              E='  +0.26%  ';  Num=$true}
    @{Row=38; D='3.50';               E='  +13.95%  '; Num=$true}
    @{Row=39; D=$null;                E='  -0.01%  ';  Num=$false}
    @{Row=40; D=$null;                E='  +0.10%  ';  Num=$false}
    @{Row=41; D=$null;                E='  +0.48%  ';  Num=$false}
    @{Row=42; D='0.999';              E='  -0.23%  ';  Num=$true}
    @{Row=43; D='2.91';               E='  -0.43%  ';  Num=$true}
    @{Row=45; D='3.49';               E='  +3.26%  ';  Num=$true}
    @{Row=46; D='3.22';               E='  +7.02%  ';  Num=$true}
    @{Row=47; D=$null;                E='  +0.74%  ';  Num=$false}
    @{Row=48; D=$null;                E='  +19.54%  '; Num=$false}
    @{Row=49; D=$null;                E='  +5.35%  ';  Num=$false}
    @{Row=50; D=$null;                E='  +0.28%  ';  Num=$false}
    @{Row=51; D="0.0${sub6}0349";     E='  -1.19%  ';  Num=$false}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($r, 4)
        if ($u.Num) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    $ws.Cells.Item($r, 5).Value = $u.E
}
